# Add the image path for the Visa "cards" document as a new value in column P
# (header "image") of the data row, mirroring the commit "Added paths to
# images in documents".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P2").Value = "C:Users/vano/Documents/GitHub/ZPI_VAF/iaff_assistant/images/Cards/VISA.png"

# Match the author's final cell selection on the new cell.
$ws.Range("P2").Select() | Out-Null
